$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 54

$ws.Cells.Item($row, 1).Value = "H8XIZ4"
$ws.Cells.Item($row, 2).Value = "Engranaje de acople de fusor de eje oscilante para HP"
$ws.Cells.Item($row, 3).Value = "Pro 400 M401 M425"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 150000
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E54-D54)*G54"
$ws.Cells.Item($row, 9).Formula = "=D54*F54"
$ws.Cells.Item($row, 10).Value = 0
